# Auto-generated Excel COM-interop script applying the workbook update
# ("Update gh-pages to output generated at e2eb17a") to before.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet 1 (展览): refresh "想去人数" (F column) interest counts ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 1243
$ws1.Cells.Item(3, 6).Value = 2006
$ws1.Cells.Item(4, 6).Value = 430
$ws1.Cells.Item(5, 6).Value = 170
$ws1.Cells.Item(6, 6).Value = 408
$ws1.Cells.Item(7, 6).Value = 45
$ws1.Cells.Item(8, 6).Value = 496
$ws1.Cells.Item(9, 6).Value = 123
$ws1.Cells.Item(10, 6).Value = 77
$ws1.Cells.Item(11, 6).Value = 155
$ws1.Cells.Item(12, 6).Value = 785
$ws1.Cells.Item(13, 6).Value = 50
$ws1.Cells.Item(15, 6).Value = 3895
$ws1.Cells.Item(16, 6).Value = 2740
$ws1.Cells.Item(17, 6).Value = 831
$ws1.Cells.Item(18, 6).Value = 606
$ws1.Cells.Item(19, 6).Value = 334
$ws1.Cells.Item(20, 6).Value = 713
$ws1.Cells.Item(21, 6).Value = 1296
$ws1.Cells.Item(22, 6).Value = 46
$ws1.Cells.Item(23, 6).Value = 675
$ws1.Cells.Item(24, 6).Value = 279
$ws1.Cells.Item(25, 6).Value = 79
$ws1.Cells.Item(26, 6).Value = 218

# ---- Sheet 2 (演出): remove the cancelled 2024.01.07 listing, shift remaining rows up ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 2).Value = '2024.01.12'
$ws2.Cells.Item(2, 3).Value = '杭州·SACG二次元live夜场（取消）'
$ws2.Cells.Item(2, 4).Value = '新北街85号三层G2-302 杭州大麦66 LIVEHOUSE'
$ws2.Cells.Item(2, 5).Value = '2024.01.12 19:00-01.12 22:00'
$ws2.Cells.Item(2, 6).Value = 42
$ws2.Cells.Item(2, 7).Value = '不可售'
$ws2.Cells.Item(2, 8).Value = $false
$ws2.Cells.Item(2, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80181&msource=Msearch_colligation'

$ws2.Cells.Item(3, 2).Value = '2024.01.13'
$ws2.Cells.Item(3, 3).Value = '杭州·《LALALAND爱乐之城》浪漫主题音乐会'
$ws2.Cells.Item(3, 4).Value = '武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）'
$ws2.Cells.Item(3, 5).Value = '2024.01.13 14:00-01.13 15:30'
$ws2.Cells.Item(3, 6).Value = 6
$ws2.Cells.Item(3, 7).Value = '90'
$ws2.Cells.Item(3, 8).Value = $false
$ws2.Cells.Item(3, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80308&msource=Msearch_colligation'

$ws2.Cells.Item(4, 2).Value = '2024.01.29'
$ws2.Cells.Item(4, 3).Value = '杭州·向着遥远的未来出发 miriちゃん生日SP'
$ws2.Cells.Item(4, 4).Value = '同协路288号 1928创意园'
$ws2.Cells.Item(4, 5).Value = '2024.01.29 18:00-01.29 21:00'
$ws2.Cells.Item(4, 6).Value = 28
$ws2.Cells.Item(4, 7).Value = '89'
$ws2.Cells.Item(4, 8).Value = $true
$ws2.Cells.Item(4, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79894&msource=Msearch_colligation'

$ws2.Cells.Item(5, 2).Value = '2024.01.31'
$ws2.Cells.Item(5, 3).Value = '杭州·《爱永恒》理查德·克莱德曼2024新年钢琴音乐会'
$ws2.Cells.Item(5, 4).Value = '武林广场29号 杭州剧院'
$ws2.Cells.Item(5, 5).Value = '2024.01.31 19:30-01.31 21:00'
$ws2.Cells.Item(5, 6).Value = 26
$ws2.Cells.Item(5, 7).Value = '不可售'
$ws2.Cells.Item(5, 8).Value = $false
$ws2.Cells.Item(5, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=76645&msource=Msearch_colligation'

$ws2.Rows.Item(6).Delete()

# ---- Sheet 3 (本地生活): refresh "想去人数" (F column) interest count ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 109

# ---- Sheet 4 (全部类型): remove the cancelled 2024.01.07 listing, shift remaining rows up, refresh counts ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 109

$ws4.Cells.Item(3, 2).Value = '2024.01.12'
$ws4.Cells.Item(3, 3).Value = '杭州·SACG二次元live夜场（取消）'
$ws4.Cells.Item(3, 4).Value = '新北街85号三层G2-302 杭州大麦66 LIVEHOUSE'
$ws4.Cells.Item(3, 5).Value = '2024.01.12 19:00-01.12 22:00'
$ws4.Cells.Item(3, 6).Value = 42
$ws4.Cells.Item(3, 7).Value = '不可售'
$ws4.Cells.Item(3, 8).Value = $false
$ws4.Cells.Item(3, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80181&msource=Msearch_colligation'

$ws4.Cells.Item(4, 2).Value = '2024.01.13'
$ws4.Cells.Item(4, 3).Value = '杭州·《LALALAND爱乐之城》浪漫主题音乐会'
$ws4.Cells.Item(4, 4).Value = '武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）'
$ws4.Cells.Item(4, 5).Value = '2024.01.13 14:00-01.13 15:30'
$ws4.Cells.Item(4, 6).Value = 6
$ws4.Cells.Item(4, 7).Value = '90'
$ws4.Cells.Item(4, 8).Value = $false
$ws4.Cells.Item(4, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80308&msource=Msearch_colligation'

$ws4.Cells.Item(5, 2).Value = '2024.01.13'
$ws4.Cells.Item(5, 3).Value = '杭州·代号鸢Only——绣衣楼过大年'
$ws4.Cells.Item(5, 4).Value = '黄姑山路51-4号 0101park'
$ws4.Cells.Item(5, 5).Value = '2024.01.13 11:00-01.13 20:00'
$ws4.Cells.Item(5, 6).Value = 1243
$ws4.Cells.Item(5, 7).Value = '98'
$ws4.Cells.Item(5, 8).Value = $false
$ws4.Cells.Item(5, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79350&msource=Msearch_colligation'

$ws4.Cells.Item(6, 2).Value = '2024.01.13'
$ws4.Cells.Item(6, 3).Value = '杭州·冒险家次元动漫展'
$ws4.Cells.Item(6, 4).Value = '沈半路171号 T-Car杭州汽车文化主题公园'
$ws4.Cells.Item(6, 5).Value = '2024.01.13 09:30-01.14 17:00'
$ws4.Cells.Item(6, 6).Value = 2006
$ws4.Cells.Item(6, 7).Value = '60'
$ws4.Cells.Item(6, 8).Value = $true
$ws4.Cells.Item(6, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79640&msource=Msearch_colligation'

$ws4.Cells.Item(7, 2).Value = '2024.01.13'
$ws4.Cells.Item(7, 3).Value = '杭州·冒险家次元动漫展 声动杭州CV专场-杜冥鸦、穆雪婷、秦紫翼、牧野冥姬'
$ws4.Cells.Item(7, 4).Value = '沈半路171号 T-Car杭州汽车文化主题公园'
$ws4.Cells.Item(7, 5).Value = '2024.01.13 10:00-01.13 16:45'
$ws4.Cells.Item(7, 6).Value = 430
$ws4.Cells.Item(7, 7).Value = '已售罄'
$ws4.Cells.Item(7, 8).Value = $true
$ws4.Cells.Item(7, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79646&msource=Msearch_colligation'

$ws4.Cells.Item(8, 2).Value = '2024.01.13'
$ws4.Cells.Item(8, 3).Value = '杭州·凡多姆海威降雪宴会·黑执事ONLY'
$ws4.Cells.Item(8, 4).Value = '凤凰御元艺术基地1138园区9号楼绿房子 迷邓花园'
$ws4.Cells.Item(8, 5).Value = '2024.01.13 12:30-01.13 21:00'
$ws4.Cells.Item(8, 6).Value = 170
$ws4.Cells.Item(8, 7).Value = '218'
$ws4.Cells.Item(8, 8).Value = $false
$ws4.Cells.Item(8, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80021&msource=Msearch_colligation'

$ws4.Cells.Item(9, 2).Value = '2024.01.13'
$ws4.Cells.Item(9, 3).Value = '杭州·温馨国乙only'
$ws4.Cells.Item(9, 4).Value = '北干街道萧杭路689号浙农东巢艺术公园 Fashion Bund时尚外滩艺术中心'
$ws4.Cells.Item(9, 5).Value = '2024.01.13 10:00-01.13 17:00'
$ws4.Cells.Item(9, 6).Value = 408
$ws4.Cells.Item(9, 7).Value = '70'
$ws4.Cells.Item(9, 8).Value = $false
$ws4.Cells.Item(9, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79169&msource=Msearch_colligation'

$ws4.Cells.Item(10, 2).Value = '2024.01.13'
$ws4.Cells.Item(10, 3).Value = '杭州·漫次元--茶话会'
$ws4.Cells.Item(10, 4).Value = '萧杭路615号2排左转到头16幢 Ciao Shed梧桐小院'
$ws4.Cells.Item(10, 5).Value = '2024.01.13 10:00-01.13 22:00'
$ws4.Cells.Item(10, 6).Value = 45
$ws4.Cells.Item(10, 7).Value = '88'
$ws4.Cells.Item(10, 8).Value = $false
$ws4.Cells.Item(10, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80448&msource=Msearch_colligation'

$ws4.Cells.Item(11, 2).Value = '2024.01.20'
$ws4.Cells.Item(11, 3).Value = '杭州·1.20新春国乙only'
$ws4.Cells.Item(11, 4).Value = '北干街道萧杭路689号浙农东巢艺术公园 Fashion Bund时尚外滩艺术中心'
$ws4.Cells.Item(11, 5).Value = '2024.01.20 10:00-01.20 17:00'
$ws4.Cells.Item(11, 6).Value = 496
$ws4.Cells.Item(11, 7).Value = '78'
$ws4.Cells.Item(11, 8).Value = $false
$ws4.Cells.Item(11, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79605&msource=Msearch_colligation'

$ws4.Cells.Item(12, 2).Value = '2024.01.20'
$ws4.Cells.Item(12, 3).Value = '杭州·动漫迷城嘉年华'
$ws4.Cells.Item(12, 4).Value = '体育场路武林广场11号 杭州大厦中央商城'
$ws4.Cells.Item(12, 5).Value = '2024.01.20 10:00-01.21 17:00'
$ws4.Cells.Item(12, 6).Value = 123
$ws4.Cells.Item(12, 7).Value = '40'
$ws4.Cells.Item(12, 8).Value = $true
$ws4.Cells.Item(12, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80216&msource=Msearch_colligation'

$ws4.Cells.Item(13, 2).Value = '2024.01.20'
$ws4.Cells.Item(13, 3).Value = '杭州·造梦探险家二次元同好会'
$ws4.Cells.Item(13, 4).Value = '振华路19号 米盒城'
$ws4.Cells.Item(13, 5).Value = '2024.01.20 10:00-01.20 17:00'
$ws4.Cells.Item(13, 6).Value = 77
$ws4.Cells.Item(13, 7).Value = '28'
$ws4.Cells.Item(13, 8).Value = $false
$ws4.Cells.Item(13, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80379&msource=Msearch_colligation'

$ws4.Cells.Item(14, 2).Value = '2024.01.21'
$ws4.Cells.Item(14, 3).Value = '杭州·春季任天堂同好会ONLY1.0'
$ws4.Cells.Item(14, 4).Value = '巧客街88号 钱塘小雅田园中心'
$ws4.Cells.Item(14, 5).Value = '2024.01.21 10:00-01.21 17:00'
$ws4.Cells.Item(14, 6).Value = 155
$ws4.Cells.Item(14, 7).Value = '69'
$ws4.Cells.Item(14, 8).Value = $false
$ws4.Cells.Item(14, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80368&msource=Msearch_colligation'

$ws4.Cells.Item(15, 2).Value = '2024.01.27'
$ws4.Cells.Item(15, 3).Value = '杭州.第32届 中二病 原神x星穹only'
$ws4.Cells.Item(15, 4).Value = '康候圣街99号 顺丰创新中心'
$ws4.Cells.Item(15, 5).Value = '2024.01.27 11:00-01.28 17:00'
$ws4.Cells.Item(15, 6).Value = 785
$ws4.Cells.Item(15, 7).Value = '50'
$ws4.Cells.Item(15, 8).Value = $true
$ws4.Cells.Item(15, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79906&msource=Msearch_colligation'

$ws4.Cells.Item(16, 2).Value = '2024.01.27'
$ws4.Cells.Item(16, 3).Value = '杭州·VOCALOID ONLY'
$ws4.Cells.Item(16, 4).Value = '九环路7号 杭州鑫牛大厦'
$ws4.Cells.Item(16, 5).Value = '2024.01.27 10:00-01.27 17:00'
$ws4.Cells.Item(16, 6).Value = 50
$ws4.Cells.Item(16, 7).Value = '78'
$ws4.Cells.Item(16, 8).Value = $true
$ws4.Cells.Item(16, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80374&msource=Msearch_colligation'

$ws4.Cells.Item(17, 2).Value = '2024.01.28'
$ws4.Cells.Item(17, 3).Value = '杭州·文豪野犬舞会ONLY:横滨晚宴（取消）'
$ws4.Cells.Item(17, 4).Value = '九龙大道227号 七里香溪别墅园'
$ws4.Cells.Item(17, 5).Value = '2024.01.28 10:00-01.28 17:00'
$ws4.Cells.Item(17, 6).Value = 851
$ws4.Cells.Item(17, 7).Value = '不可售'
$ws4.Cells.Item(17, 8).Value = $true
$ws4.Cells.Item(17, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=77627&msource=Msearch_colligation'

$ws4.Cells.Item(18, 2).Value = '2024.01.29'
$ws4.Cells.Item(18, 3).Value = '杭州·向着遥远的未来出发 miriちゃん生日SP'
$ws4.Cells.Item(18, 4).Value = '同协路288号 1928创意园'
$ws4.Cells.Item(18, 5).Value = '2024.01.29 18:00-01.29 21:00'
$ws4.Cells.Item(18, 6).Value = 28
$ws4.Cells.Item(18, 7).Value = '89'
$ws4.Cells.Item(18, 8).Value = $true
$ws4.Cells.Item(18, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79894&msource=Msearch_colligation'

$ws4.Cells.Item(19, 2).Value = '2024.01.31'
$ws4.Cells.Item(19, 3).Value = '杭州·《爱永恒》理查德·克莱德曼2024新年钢琴音乐会'
$ws4.Cells.Item(19, 4).Value = '武林广场29号 杭州剧院'
$ws4.Cells.Item(19, 5).Value = '2024.01.31 19:30-01.31 21:00'
$ws4.Cells.Item(19, 6).Value = 26
$ws4.Cells.Item(19, 7).Value = '不可售'
$ws4.Cells.Item(19, 8).Value = $false
$ws4.Cells.Item(19, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=76645&msource=Msearch_colligation'

$ws4.Cells.Item(20, 2).Value = '2024.02.03'
$ws4.Cells.Item(20, 3).Value = '杭州·AP动漫游戏嘉年华'
$ws4.Cells.Item(20, 4).Value = '飞虹路3号 杭州奥体中心综合训练馆'
$ws4.Cells.Item(20, 5).Value = '2024.02.03 09:00-02.04 17:00'
$ws4.Cells.Item(20, 6).Value = 3895
$ws4.Cells.Item(20, 7).Value = '70'
$ws4.Cells.Item(20, 8).Value = $true
$ws4.Cells.Item(20, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79837&msource=Msearch_colligation'

$ws4.Cells.Item(21, 2).Value = '2024.02.03'
$ws4.Cells.Item(21, 3).Value = '杭州·樱之弦世界动漫游戏博览会（取消）'
$ws4.Cells.Item(21, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws4.Cells.Item(21, 5).Value = '2024.02.03 10:00-02.04 17:00'
$ws4.Cells.Item(21, 6).Value = 2740
$ws4.Cells.Item(21, 7).Value = '不可售'
$ws4.Cells.Item(21, 8).Value = $false
$ws4.Cells.Item(21, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=78232&msource=Msearch_colligation'

$ws4.Cells.Item(22, 2).Value = '2024.02.03'
$ws4.Cells.Item(22, 3).Value = '杭州·浙江蔚蓝档案only'
$ws4.Cells.Item(22, 4).Value = '石祥路242号 首开公园'
$ws4.Cells.Item(22, 5).Value = '2024.02.03 10:00-02.03 18:00'
$ws4.Cells.Item(22, 6).Value = 831
$ws4.Cells.Item(22, 7).Value = '已售罄'
$ws4.Cells.Item(22, 8).Value = $false
$ws4.Cells.Item(22, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=77937&msource=Msearch_colligation'

$ws4.Cells.Item(23, 2).Value = '2024.02.04'
$ws4.Cells.Item(23, 3).Value = '杭州·偶像梦幻祭ONLY'
$ws4.Cells.Item(23, 4).Value = '石祥路242号 首开公园'
$ws4.Cells.Item(23, 5).Value = '2024.02.04 10:00-02.04 18:00'
$ws4.Cells.Item(23, 6).Value = 606
$ws4.Cells.Item(23, 7).Value = '68'
$ws4.Cells.Item(23, 8).Value = $false
$ws4.Cells.Item(23, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=77815&msource=Msearch_colligation'

$ws4.Cells.Item(24, 2).Value = '2024.02.16'
$ws4.Cells.Item(24, 3).Value = '杭州·PJSK only展'
$ws4.Cells.Item(24, 4).Value = '同协路288号 1928创意园'
$ws4.Cells.Item(24, 5).Value = '2024.02.16 10:00-02.16 18:00'
$ws4.Cells.Item(24, 6).Value = 334
$ws4.Cells.Item(24, 7).Value = '已售罄'
$ws4.Cells.Item(24, 8).Value = $false
$ws4.Cells.Item(24, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80124&msource=Msearch_colligation'

$ws4.Cells.Item(25, 2).Value = '2024.02.16'
$ws4.Cells.Item(25, 3).Value = '杭州·第34届中二病动漫游戏展'
$ws4.Cells.Item(25, 4).Value = '康候圣街99号 顺丰创新中心'
$ws4.Cells.Item(25, 5).Value = '2024.02.16 11:00-02.16 17:00'
$ws4.Cells.Item(25, 6).Value = 713
$ws4.Cells.Item(25, 7).Value = '60'
$ws4.Cells.Item(25, 8).Value = $false
$ws4.Cells.Item(25, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79971&msource=Msearch_colligation'

$ws4.Cells.Item(26, 2).Value = '2024.02.17'
$ws4.Cells.Item(26, 3).Value = '杭州·6th YH樱花动漫游戏文化节'
$ws4.Cells.Item(26, 4).Value = '亚太路湘湖3期东南侧约290米 原创壹号羽毛球馆'
$ws4.Cells.Item(26, 5).Value = '2024.02.17 10:00-02.18 17:00'
$ws4.Cells.Item(26, 6).Value = 1296
$ws4.Cells.Item(26, 7).Value = '65'
$ws4.Cells.Item(26, 8).Value = $true
$ws4.Cells.Item(26, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80324&msource=Msearch_colligation'

$ws4.Cells.Item(27, 2).Value = '2024.02.24'
$ws4.Cells.Item(27, 3).Value = '杭州·次元幻想动漫游戏嘉年华'
$ws4.Cells.Item(27, 4).Value = '德胜东路2539号 梦马汽车小镇'
$ws4.Cells.Item(27, 5).Value = '2024.02.24 10:00-02.25 17:00'
$ws4.Cells.Item(27, 6).Value = 46
$ws4.Cells.Item(27, 7).Value = '58'
$ws4.Cells.Item(27, 8).Value = $false
$ws4.Cells.Item(27, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80425&msource=Msearch_colligation'

$ws4.Cells.Item(28, 2).Value = '2024.02.24'
$ws4.Cells.Item(28, 3).Value = '杭州第35届 中二病 原神x星穹only'
$ws4.Cells.Item(28, 4).Value = '康候圣街99号 顺丰创新中心'
$ws4.Cells.Item(28, 5).Value = '2024.02.24 11:00-02.24 17:00'
$ws4.Cells.Item(28, 6).Value = 675
$ws4.Cells.Item(28, 7).Value = '50'
$ws4.Cells.Item(28, 8).Value = $true
$ws4.Cells.Item(28, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79890&msource=Msearch_colligation'

$ws4.Cells.Item(29, 2).Value = '2024.03.02'
$ws4.Cells.Item(29, 3).Value = '杭州·杭州灵能百分百only'
$ws4.Cells.Item(29, 4).Value = '水博大道118号 宝盛水博园大酒店'
$ws4.Cells.Item(29, 5).Value = '2024.03.02 09:30-03.02 17:00'
$ws4.Cells.Item(29, 6).Value = 279
$ws4.Cells.Item(29, 7).Value = '89'
$ws4.Cells.Item(29, 8).Value = $false
$ws4.Cells.Item(29, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=80099&msource=Msearch_colligation'

$ws4.Cells.Item(30, 2).Value = '2024.03.09'
$ws4.Cells.Item(30, 3).Value = '杭州·异次结界动漫嘉年华'
$ws4.Cells.Item(30, 4).Value = '德胜东路2539号 梦马汽车小镇'
$ws4.Cells.Item(30, 5).Value = '2024.03.09 10:00-03.10 17:00'
$ws4.Cells.Item(30, 6).Value = 79
$ws4.Cells.Item(30, 7).Value = '65'
$ws4.Cells.Item(30, 8).Value = $false
$ws4.Cells.Item(30, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=79935&msource=Msearch_colligation'

$ws4.Cells.Item(31, 2).Value = '2024.04.04'
$ws4.Cells.Item(31, 3).Value = '杭州·第九届萌次元动漫嘉年华'
$ws4.Cells.Item(31, 4).Value = '长乐路29号五组2幢 杭州运河文化发布中心'
$ws4.Cells.Item(31, 5).Value = '2024.04.04 10:00-04.05 17:00'
$ws4.Cells.Item(31, 6).Value = 218
$ws4.Cells.Item(31, 7).Value = '不可售'
$ws4.Cells.Item(31, 8).Value = $false
$ws4.Cells.Item(31, 9).Value = 'https://show.bilibili.com/platform/detail.html?id=78866&msource=Msearch_colligation'

$ws4.Rows.Item(32).Delete()
